# Fruta / hortaliza, semanal
# Insert three new weekly price rows (Femacal de La Calera - Arándano (blue))
# at sheet rows 118-120, pushing the existing rows 118-140 down to 121-143.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the current row 118 (each Insert() call pushes
# everything at/below that row index down by one, so calling it three times
# at the same target row opens up a 3-row gap at 118-120).
$ws.Rows.Item(118).Insert()
$ws.Rows.Item(118).Insert()
$ws.Rows.Item(118).Insert()

# New row 118
$ws.Cells.Item(118, 1).Value = 3
$ws.Cells.Item(118, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44504
$ws.Cells.Item(118, 5).Value = 5
$ws.Cells.Item(118, 6).Value = "Fruta"
$ws.Cells.Item(118, 7).Value = 100101
$ws.Cells.Item(118, 8).Value = "Berries"
$ws.Cells.Item(118, 9).Value = 100101001
$ws.Cells.Item(118, 10).Value = "Arándano (blue)"
$ws.Cells.Item(118, 11).Value = "Sin especificar"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 50
$ws.Cells.Item(118, 14).Value = 10000
$ws.Cells.Item(118, 15).Value = 10000
$ws.Cells.Item(118, 16).Value = 10000
$ws.Cells.Item(118, 17).Value = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(118, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(118, 19).Value = 6667
$ws.Cells.Item(118, 20).Value = 1.5

# New row 119
$ws.Cells.Item(119, 1).Value = 3
$ws.Cells.Item(119, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44504
$ws.Cells.Item(119, 5).Value = 5
$ws.Cells.Item(119, 6).Value = "Fruta"
$ws.Cells.Item(119, 7).Value = 100101
$ws.Cells.Item(119, 8).Value = "Berries"
$ws.Cells.Item(119, 9).Value = 100101001
$ws.Cells.Item(119, 10).Value = "Arándano (blue)"
$ws.Cells.Item(119, 11).Value = "Sin especificar"
$ws.Cells.Item(119, 12).Value = "Primera"
$ws.Cells.Item(119, 13).Value = 56
$ws.Cells.Item(119, 14).Value = 10000
$ws.Cells.Item(119, 15).Value = 10000
$ws.Cells.Item(119, 16).Value = 10000
$ws.Cells.Item(119, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(119, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(119, 19).Value = 5000
$ws.Cells.Item(119, 20).Value = 2

# New row 120
$ws.Cells.Item(120, 1).Value = 3
$ws.Cells.Item(120, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44504
$ws.Cells.Item(120, 5).Value = 5
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100101
$ws.Cells.Item(120, 8).Value = "Berries"
$ws.Cells.Item(120, 9).Value = 100101001
$ws.Cells.Item(120, 10).Value = "Arándano (blue)"
$ws.Cells.Item(120, 11).Value = "Sin especificar"
$ws.Cells.Item(120, 12).Value = "Segunda"
$ws.Cells.Item(120, 13).Value = 50
$ws.Cells.Item(120, 14).Value = 8000
$ws.Cells.Item(120, 15).Value = 8000
$ws.Cells.Item(120, 16).Value = 8000
$ws.Cells.Item(120, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(120, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(120, 19).Value = 4000
$ws.Cells.Item(120, 20).Value = 2
